# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Swap the display order of "Arroyo de la Luz" (row 56) and "La Palma" (row 57)
# so that "La Palma" now appears before "Arroyo de la Luz" in the shared
# string table (their row/numeric data stays the same).
$ws.Range("A56").Value = "La Palma"
$ws.Range("A57").Value = "Arroyo de la Luz"

# Update the "Datos actualizados" timestamp string in cell A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 18:16"

# Update Navarra row (row 8): Recuperados 585 -> 583, Muertes 6 -> 8
$ws.Range("D8").Value = 583
$ws.Range("E8").Value = 8
